$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 41281
$ws.Range("D2").Value = 59618478
$ws.Range("C3").Value = 98106
$ws.Range("D3").Value = 143728495
$ws.Range("C4").Value = 33384
$ws.Range("D4").Value = 49415697
$ws.Range("C5").Value = 9511
$ws.Range("D5").Value = 14126742
$ws.Range("C6").Value = 2326
$ws.Range("D6").Value = 3455920
$ws.Range("C7").Value = 234
$ws.Range("D7").Value = 346093
$ws.Range("C12").Value = 44450
$ws.Range("D12").Value = 60183064
$ws.Range("C13").Value = 10455
$ws.Range("D13").Value = 15101269
$ws.Range("C14").Value = 27656
$ws.Range("D14").Value = 40526208
$ws.Range("C15").Value = 8775
$ws.Range("D15").Value = 13021422
$ws.Range("C16").Value = 2320
$ws.Range("D16").Value = 3447603
$ws.Range("C17").Value = 469
$ws.Range("D17").Value = 692623
$ws.Range("C20").Value = 10914
$ws.Range("D20").Value = 14378386
$ws.Range("C21").Value = 14382
$ws.Range("D21").Value = 20737319
$ws.Range("C22").Value = 33531
$ws.Range("D22").Value = 49161306
$ws.Range("C23").Value = 10792
$ws.Range("D23").Value = 16036229
$ws.Range("C24").Value = 2857
$ws.Range("D24").Value = 4245615
$ws.Range("C25").Value = 587
$ws.Range("D25").Value = 874592
$ws.Range("C27").Value = 12433
$ws.Range("D27").Value = 16537066
$ws.Range("C28").Value = 8380
$ws.Range("D28").Value = 12118361
$ws.Range("C29").Value = 24080
$ws.Range("D29").Value = 35327438
$ws.Range("C30").Value = 8280
$ws.Range("D30").Value = 12308707
$ws.Range("C31").Value = 2114
$ws.Range("D31").Value = 3153208
$ws.Range("C32").Value = 417
$ws.Range("D32").Value = 616249
$ws.Range("C34").Value = 8934
$ws.Range("D34").Value = 11777181
$ws.Range("C35").Value = 3639
$ws.Range("D35").Value = 5254692
$ws.Range("C36").Value = 8559
$ws.Range("D36").Value = 12504993
$ws.Range("C37").Value = 3393
$ws.Range("D37").Value = 5031738
$ws.Range("C38").Value = 867
$ws.Range("D38").Value = 1291555
$ws.Range("C41").Value = 2722
$ws.Range("D41").Value = 3666549
$ws.Range("C42").Value = 18693
$ws.Range("D42").Value = 26994773
$ws.Range("C43").Value = 54515
$ws.Range("D43").Value = 79878702
$ws.Range("C44").Value = 19993
$ws.Range("D44").Value = 29681697
$ws.Range("C45").Value = 6000
$ws.Range("D45").Value = 8923845
$ws.Range("C46").Value = 1399
$ws.Range("D46").Value = 2088144
$ws.Range("C50").Value = 18114
$ws.Range("D50").Value = 24001097
$ws.Range("C51").Value = 2352
$ws.Range("D51").Value = 3412031
$ws.Range("C52").Value = 7920
$ws.Range("D52").Value = 11630332
$ws.Range("C53").Value = 2654
$ws.Range("D53").Value = 3961133
$ws.Range("C54").Value = 831
$ws.Range("D54").Value = 1241415
$ws.Range("C55").Value = 216
$ws.Range("D55").Value = 319448
$ws.Range("C57").Value = 7907
$ws.Range("D57").Value = 10883295
$ws.Range("C58").Value = 1623
$ws.Range("D58").Value = 3239318
$ws.Range("C59").Value = 3873
$ws.Range("D59").Value = 7692524
$ws.Range("C60").Value = 1531
$ws.Range("D60").Value = 3049462
$ws.Range("C61").Value = 508
$ws.Range("D61").Value = 1005083
$ws.Range("C62").Value = 183
$ws.Range("D62").Value = 379487
$ws.Range("C64").Value = 2521
$ws.Range("D64").Value = 4652238
$ws.Range("C65").Value = 16840
$ws.Range("D65").Value = 24302084
$ws.Range("C66").Value = 48026
$ws.Range("D66").Value = 70208252
$ws.Range("C67").Value = 16770
$ws.Range("D67").Value = 24914976
$ws.Range("C68").Value = 4888
$ws.Range("D68").Value = 7279524
$ws.Range("C69").Value = 1071
$ws.Range("D69").Value = 1592199
$ws.Range("C70").Value = 97
$ws.Range("D70").Value = 142830
$ws.Range("C73").Value = 16064
$ws.Range("D73").Value = 21084151
$ws.Range("C74").Value = 61704
$ws.Range("D74").Value = 89694921
$ws.Range("C75").Value = 168334
$ws.Range("D75").Value = 247748657
$ws.Range("C76").Value = 71703
$ws.Range("D76").Value = 106789531
$ws.Range("C77").Value = 23368
$ws.Range("D77").Value = 34907040
$ws.Range("C78").Value = 5981
$ws.Range("D78").Value = 8930544
$ws.Range("C79").Value = 426
$ws.Range("D79").Value = 633460
$ws.Range("C84").Value = 6
$ws.Range("D84").Value = 9000
$ws.Range("C85").Value = 60733
$ws.Range("D85").Value = 81998609
$ws.Range("C86").Value = 5111
$ws.Range("D86").Value = 7403991
$ws.Range("C87").Value = 12574
$ws.Range("D87").Value = 18464765
$ws.Range("C88").Value = 4123
$ws.Range("D88").Value = 6143640
$ws.Range("C89").Value = 1449
$ws.Range("D89").Value = 2164611
$ws.Range("C90").Value = 341
$ws.Range("D90").Value = 508012
$ws.Range("C93").Value = 5849
$ws.Range("D93").Value = 7842511
$ws.Range("C94").Value = 1813
$ws.Range("D94").Value = 2611771
$ws.Range("C95").Value = 5788
$ws.Range("D95").Value = 8531424
$ws.Range("C96").Value = 2088
$ws.Range("D96").Value = 3107781
$ws.Range("C97").Value = 767
$ws.Range("D97").Value = 1149460
$ws.Range("C98").Value = 214
$ws.Range("D98").Value = 323113
$ws.Range("C101").Value = 3940
$ws.Range("D101").Value = 5225260
$ws.Range("C102").Value = 903
$ws.Range("D102").Value = 1731061
$ws.Range("C103").Value = 598
$ws.Range("D103").Value = 1211142
$ws.Range("C104").Value = 217
$ws.Range("D104").Value = 432194
$ws.Range("C105").Value = 70
$ws.Range("D105").Value = 136500
$ws.Range("C107").Value = 11831
$ws.Range("D107").Value = 17147820
$ws.Range("C108").Value = 31006
$ws.Range("D108").Value = 45511349
$ws.Range("C109").Value = 10374
$ws.Range("D109").Value = 15424749
$ws.Range("C110").Value = 2873
$ws.Range("D110").Value = 4282571
$ws.Range("C111").Value = 556
$ws.Range("D111").Value = 827953
$ws.Range("C115").Value = 10450
$ws.Range("D115").Value = 13753668
$ws.Range("C116").Value = 33189
$ws.Range("D116").Value = 47819161
$ws.Range("C117").Value = 70576
$ws.Range("D117").Value = 103230034
$ws.Range("C118").Value = 22632
$ws.Range("D118").Value = 33617217
$ws.Range("C119").Value = 6484
$ws.Range("D119").Value = 9653232
$ws.Range("C120").Value = 1280
$ws.Range("D120").Value = 1911737
$ws.Range("C121").Value = 118
$ws.Range("D121").Value = 173395
$ws.Range("C125").Value = 27575
$ws.Range("D125").Value = 36725264
$ws.Range("C126").Value = 39626
$ws.Range("D126").Value = 57127525
$ws.Range("C127").Value = 82837
$ws.Range("D127").Value = 121044687
$ws.Range("C128").Value = 25387
$ws.Range("D128").Value = 37670424
$ws.Range("C129").Value = 6921
$ws.Range("D129").Value = 10284009
$ws.Range("C130").Value = 1454
$ws.Range("D130").Value = 2154096
$ws.Range("C131").Value = 84
$ws.Range("D131").Value = 124228
$ws.Range("C134").Value = 34041
$ws.Range("D134").Value = 45085585
$ws.Range("C135").Value = 14440
$ws.Range("D135").Value = 20892743
$ws.Range("C136").Value = 34429
$ws.Range("D136").Value = 50535601
$ws.Range("C137").Value = 12147
$ws.Range("D137").Value = 18047073
$ws.Range("C138").Value = 3220
$ws.Range("D138").Value = 4799875
$ws.Range("C139").Value = 581
$ws.Range("D139").Value = 865490
$ws.Range("C143").Value = 11547
$ws.Range("D143").Value = 15349993
$ws.Range("C144").Value = 38835
$ws.Range("D144").Value = 56073554
$ws.Range("C145").Value = 88876
$ws.Range("D145").Value = 130128623
$ws.Range("C146").Value = 26441
$ws.Range("D146").Value = 39264843
$ws.Range("C147").Value = 7028
$ws.Range("D147").Value = 10472072
$ws.Range("C148").Value = 1649
$ws.Range("D148").Value = 2449143
$ws.Range("C149").Value = 108
$ws.Range("D149").Value = 161630
$ws.Range("C151").Value = 31511
$ws.Range("D151").Value = 42368595

Write-Host "Updated 112 rows"